$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the policy number used by the "preproduccion" row (row 2) ---
# Leading apostrophe keeps it text (same quotePrefix style as before) instead
# of Excel re-interpreting it as a number and dropping the leading zero.
$ws.Range("E2").Value = "'04104013002"

# --- Stash B2's current formatting (its "Hipervinculo" style) so we can
# restore it after the hyperlink collection gets rebuilt below ---
$ws.Range("Z1").Value = "fmt"
$ws.Range("B2").Copy($ws.Range("Z1"))

# --- Remove the whole second data row (the "Oci2 / suraqa" environment) ---
$ws.Rows("4").Delete()

# Rebuilding the single remaining hyperlink: this runtime's Hyperlinks.Delete
# clears every hyperlink on the sheet rather than just the targeted range, so
# clear them all and re-add the one that must survive (B2).
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do") | Out-Null

# Restore B2's original cell format/style (re-adding the hyperlink resets it)
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Z1").Clear()

# --- Match the saved selection state ---
$ws.Range("G10").Select()
